$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header row ---
# Column A/B header labels changed from Controls/Experimentals to WT/MSxc
$ws.Range("A1").Value = "WT"
$ws.Range("B1").Value = "MSxc"
$ws.Range("C1").Value = "Acceptable MSNs"
# New column D header - flags programs to exclude from parsing
$ws.Range("D1").Value = "Do Not Parse for Test Breakdown"
# The old D1 note now moves to E1
$ws.Range("E1").Value = "<--- dictate what programs collect data in a similar way"

# --- Data rows: WT / MSxc mouse IDs and their acceptable MSN program names ---
$ws.Range("A2").Value = "M2-4624"
$ws.Range("B2").Value = "M2-4645"
$ws.Range("C2").Value = "25P STOP LT 3 Greg ZD v3"
$ws.Range("D2").Value = "Lever Training 2 Greg"

$ws.Range("A3").Value = "M2-4670"
$ws.Range("B3").Value = "M2-4626"
$ws.Range("C3").Value = "25P STOP LT 3 Greg ZD BoxA2"
$ws.Range("D3").Value = "Lever Training 2 Greg Box A2"

$ws.Range("A4").Value = "M2-4780"
$ws.Range("B4").Value = "M2-4600"
$ws.Range("C4").Value = "25P STOP 200ms TEST Greg SSRT BoxA4 v2"

$ws.Range("A5").Value = "M2-4795"
$ws.Range("B5").Value = "M2-4599"
$ws.Range("C5").Value = "25P STOP 50ms TEST Greg SSRT BoxA4 v2"

$ws.Range("A6").Value = "M2-4860"
$ws.Range("B6").Value = "M2-4874"
$ws.Range("C6").Value = "25P STOP ZD TEST Greg SSRT BoxA4 v2"

$ws.Range("A7").Value = "M2-4911"
$ws.Range("B7").Value = "M2-4887"
$ws.Range("C7").Value = "25P STOP Baseline Greg SSRT BoxA4 v2"

$ws.Range("A8").Value = "M2-4861"
$ws.Range("B8").Value = "M2-4886"
$ws.Range("C8").Value = "25P STOP LT 3 Greg ZD v3 BoxA4"

$ws.Range("A9").Value = "M2-4864"
$ws.Range("B9").Value = "M2-4917"
$ws.Range("C9").Value = "Lever Training 2 Greg"

$ws.Range("A10").Value = "M2-4975"
$ws.Range("C10").Value = "25P STOP 100ms TEST Greg SSRT BoxA2"

# --- Remaining "Acceptable MSNs" list entries (column C only) ---
$ws.Range("C11").Value = "25P STOP 300ms TEST Greg SSRT BoxA2"
$ws.Range("C12").Value = "25P STOP Baseline Greg SSRT BoxA2 v2"
$ws.Range("C13").Value = "Lever Training 2 Greg Box A2"
$ws.Range("C14").Value = "Lever Training 1"
$ws.Range("C15").Value = "LT 3 Greg ZD"
$ws.Range("C16").Value = "LT 3 Greg ZD BoxA2"
$ws.Range("C17").Value = "LT 3 Greg ZD BoxA4"
$ws.Range("C18").Value = "Stop lever delay LT 3 Greg ZD v3"
$ws.Range("C19").Value = "Stop lever delay LT 3 Greg ZD v3 A4"
$ws.Range("C20").Value = "Stop lever delay LT 3 Greg ZD v3 A2"
$ws.Range("C21").Value = "30sec pun LT 3 Greg ZD v3 BoxA2"
$ws.Range("C22").Value = "30sec pun LT 3 Greg ZD v3"
$ws.Range("C23").Value = "30sec pun LT 3 Greg ZD v3 BoxA4"
$ws.Range("C24").Value = "LT 3 ZD No Stop Lever BoxA4"
$ws.Range("C25").Value = "LT 3 ZD No Stop Lever"
$ws.Range("C26").Value = "LT 3 ZD No Stop Lever BoxA2"
$ws.Range("C27").Value = "Stop Lever Delay LT3 v1 Box A2"
$ws.Range("C28").Value = "Stop Lever Delay LT3 v1"
$ws.Range("C29").Value = "4911 Baseline SSRT BoxA2 v3"

# --- Formatting: header row bold (matches existing A1:C1 style) & new column widths ---
$ws.Range("A1:E1").Font.Bold = $true
$ws.Columns("C:D").EntireColumn.AutoFit()

# --- Selection matches the author's saved cursor position ---
$ws.Range("D11").Select()
